$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from the dated "18-01-2025" tab to the generic "Total"
# tab. Renaming through the COM object model keeps every workbook-scoped
# defined name (hidden AutoFilter/wvu.PrintTitles bookmarks, etc.) pointed
# at the sheet automatically.
$ws.Name = "Total"

# The rename above does not reach the two "special" defined names that back
# Print_Area / Print_Titles (they are driven off PageSetup, not the raw
# Names collection), so refresh those explicitly to repoint them at the
# renamed sheet with the same ranges they already had.
$ws.PageSetup.PrintArea = "A1:M9"
$ws.PageSetup.PrintTitleRows = "$8:$9"

